$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "592.84")
# are stored verbatim as text, matching the source workbook's inlineStr cells,
# instead of being auto-coerced to floating point numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.929.65"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.510.70"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "592.84"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "173.14"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("E9").Value = "  +6.37%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "0.436"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "4.122.96"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "28.92"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "66.963.32"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "3.418.93"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "6.31"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "14.22"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "395.31"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "7.97"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "73.23"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "0.539"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "0.0000121"
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("D26").Value = "10.18"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "6.25"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "23.90"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "7.38"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").Value = "163.07"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "0.893"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "6.88"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").Value = "27.81"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").Value = "4.68"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "0.0742"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "26.39"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "2.802.83"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "42.84"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "0.0305"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").Value = "338.64"
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "33.40"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "6.51"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.845"
$ws.Range("E51").Value = "  -1.10%  "

# Restore the default cell style so no stray number-format style lingers
# on cells that did not have one in the source file.
$dRange.Style = "Normal"
